$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the literal string (e.g. "63.20") is preserved,
# then restore the Normal style so no stray formatting is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D5").Value = "524.64"
$ws.Range("D6").Value = "147.19"
$ws.Range("D16").Value = "21.39"
$ws.Range("D19").Value = "352.27"
$ws.Range("D20").Value = "4.55"
$ws.Range("D22").Value = "6.36"
$ws.Range("D24").Value = "63.20"
$ws.Range("D25").Value = "0.423"
$ws.Range("D27").Value = "0.993"
$ws.Range("D29").Value = "7.36"
$ws.Range("D30").Value = "6.84"
$ws.Range("D33").Value = "19.16"
$ws.Range("D35").Value = "4.31"
$ws.Range("D37").Value = "0.957"
$ws.Range("D38").Value = "0.880"
$ws.Range("D39").Value = "1.53"
$ws.Range("D40").Value = "36.92"
$ws.Range("D41").Value = "3.72"
$ws.Range("D42").Value = "285.59"
$ws.Range("D43").Value = "20.14"
$ws.Range("D44").Value = "0.0994"
$ws.Range("D45").Value = "0.613"
$ws.Range("D48").Value = "4.94"
$ws.Range("D49").Value = "0.0542"
$ws.Range("D51").Value = "19.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# Remaining cells: plain text / already non-numeric-looking strings.
$ws.Range("D2").Value = "60.313.06"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "2.688.26"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "2.709.85"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "3.164.29"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "60.325.57"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "2.696.12"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0822"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("E30").Value = "  +5.67%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  +5.07%  "
$ws.Range("E36").Value = "  +8.01%  "
$ws.Range("E37").Value = "  -7.05%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E39").Value = "  +7.33%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "2.133.97"
$ws.Range("E47").Value = "  +5.55%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E51").Value = "  +4.63%  "
